$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Population size (P): 4190 -> 50000
$ws.Range("F3").Value = 50000

# Remove "Degree of Confidence (Z):" label and its value; row becomes blank
$ws.Range("C4").ClearContents()
$ws.Range("F4").ClearContents()

# Expected Occurrence (p): 0.05 -> 0.0759
$ws.Range("F5").Value = 0.075899999999999995

# Size of sample (n): 45 -> 500
$ws.Range("E20").Value = 500

# # matching criteria: 2 -> 100
$ws.Range("E21").Value = 100

$ws.Calculate()
